# Insert a new weekly price record for "Femacal de La Calera - Achicoria"
# at row 153, pushing the previously existing rows 153..203 down to 154..204.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 153 (default shifts the existing rows down).
$ws.Rows.Item(153).Insert()

# Fill in the new row with the new data record.
$ws.Cells.Item(153, 1).Value  = 3
$ws.Cells.Item(153, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(153, 3).Value  = "Coquimbo"
$ws.Cells.Item(153, 4).Value  = 44627
$ws.Cells.Item(153, 5).Value  = 5
$ws.Cells.Item(153, 6).Value  = 100112010
$ws.Cells.Item(153, 7).Value  = "Achicoria"
$ws.Cells.Item(153, 8).Value  = "Sin especificar"
$ws.Cells.Item(153, 9).Value  = "Primera"
$ws.Cells.Item(153, 10).Value = 81
$ws.Cells.Item(153, 11).Value = 7000
$ws.Cells.Item(153, 12).Value = 7500
$ws.Cells.Item(153, 13).Value = 7222
$ws.Cells.Item(153, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(153, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(153, 16).Value = 451
$ws.Cells.Item(153, 17).Value = 16
$ws.Cells.Item(153, 18).Value = "Hortaliza"
